# "Setup of Bitmap Ops" — insert a new Text_GetRect entry into the Text_*
# function list in column A (rows 8-17), add the "x" marker cells in
# column B for the newly-aligned rows, drop the stray A28 border-only
# cell, and move the active selection to B10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the Text_* names in column A down by one starting at row 10 ---
# (only column A moves; the neighbouring C/E/G/I/K/M/O columns stay put)
# Walk bottom-up so we never clobber a value before it has been copied.
$ws.Range("A17").Value = $ws.Range("A16").Value()
$ws.Range("A16").Value = $ws.Range("A15").Value()
$ws.Range("A15").Value = $ws.Range("A14").Value()
$ws.Range("A14").Value = $ws.Range("A13").Value()
$ws.Range("A13").Value = $ws.Range("A12").Value()
$ws.Range("A12").Value = $ws.Range("A11").Value()
$ws.Range("A11").Value = $ws.Range("A10").Value()
$ws.Range("A10").Value = "Text_GetRect"

# --- Add the "x" marker cells in column B, copying B1's look ("x", centred/bold/shaded) ---
$ws.Range("B1").Copy()
$ws.Range("B8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteAll)
$ws.Range("B9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteAll)
$ws.Range("B11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteAll)
$ws.Range("B12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteAll)
$ws.Range("B13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteAll)
$ws.Range("B14").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteAll)
$ws.Range("B15").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteAll)
$ws.Range("B16").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteAll)

# --- Drop the now-stray A28 border-only cell ---
$ws.Range("A28").Clear()

# --- Move the active selection ---
$ws.Range("B10").Select()
